$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("parameters")
$ws2 = $wb.Worksheets.Item("bias")
$ws5 = $wb.Worksheets.Item("climate")

# ---------------------------------------------------------------------------
# Move the 17 "bias" parameter rows (rows 2-18 of the `bias` sheet) up into
# the `parameters` sheet, right after the existing rows, starting at row 63.
# This pushes the former rows 63-66 of `parameters` down to rows 80-83.
# ---------------------------------------------------------------------------

# 1) Make room in "parameters" for the 17 rows being relocated.
$ws1.Range("A63:F79").Insert(-4121)

# 2) Copy the values from "bias" into the freshly inserted rows.
$ws2.Range("A2:F18").Copy($ws1.Range("A63:F79"))

# Two of those source rows had no value in column D - remove the stray
# empty cells Copy leaves behind so the rows match exactly.
$ws1.Range("D74").ClearContents()
$ws1.Range("D79").ClearContents()

# 3) Remove the now-duplicated rows from "bias", shifting the remaining
#    rows (formerly 19-48) up so they become rows 2-31.
$ws2.Range("A2:F18").Delete(-4162)

# ---------------------------------------------------------------------------
# View / selection state.
# ---------------------------------------------------------------------------

# "climate" was the active sheet before; update its scroll position and
# selection, then let it lose the active-tab status below.
$ws5.Activate()
$excel.ActiveWindow.ScrollRow = 93
$excel.ActiveWindow.ScrollColumn = 1
$ws5.Range("M145").Select()

# "parameters" keeps its own scroll position / selection.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A83").Select()

# "bias" becomes the active / selected sheet and cell.
$ws2.Activate()
$ws2.Range("E37").Select()
